$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data window shifts from 2000-2020 (rows 2-22) to 2010-2022 (rows 2-14):
#  - rows 15-22 (old 2013-2020 data) are no longer needed as separate rows
#    since the whole block is being replaced, so first drop the extra rows
#    at the bottom (15-22) to shrink the sheet to 14 rows.
#  - then overwrite rows 2-14 with the updated 2010-2022 figures.

$ws.Rows("15:22").Delete()

$data = @(
    @("2010年", 100, 9.3251769405, 44.1769919214, 46.4978311381),
    @("2011年", 100, 9.177653795399999, 44.2930567346, 46.52928947),
    @("2012年", 100, 9.1137138703, 45.4633053876, 45.4229807421),
    @("2013年", 100, 8.942893962299999, 46.8804015014, 44.1767045363),
    @("2014年", 100, 8.6434913962, 48.2709404943, 43.0855681094),
    @("2015年", 100, 8.387014882100001, 50.7716451484, 40.8413399694),
    @("2016年", 100, 8.0572875347, 52.3620909629, 39.5806215024),
    @("2017年", 100, 7.4635650166, 52.6847365342, 39.8516984492),
    @("2018年", 100, 7.0430202532, 53.2699679234, 39.6870118234),
    @("2019年", 100, 7.14369027759977, 54.268904266364, 38.5874054560362),
    @("2020年", 100, 7.69864253413716, 54.4585357459341, 37.8428217199287),
    @("2021年", 100, 7.2, 53.5, 39.3),
    @("2022年", 100, 7.3, 52.8, 39.9)
)

$r = 2
foreach ($row in $data) {
    $ws.Range("A$r").Value2 = $row[0]
    $ws.Range("B$r").Value2 = $row[1]
    $ws.Range("C$r").Value2 = $row[2]
    $ws.Range("D$r").Value2 = $row[3]
    $ws.Range("E$r").Value2 = $row[4]
    $r = $r + 1
}
